# Add the new event row (row 11) to the "all_Events" sheet.
# Source data (my_Events.xlsx) stores every cell as text, so the numeric-
# looking values (the big numeric ID, day/month/year, reminder days) must be
# written as text too -- otherwise Excel would coerce them to numbers and the
# 19-digit ID would lose precision via scientific notation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_Events")

$row = 11
$rng = $ws.Range("A11:F11")

# Force text storage for the whole new row so numeric-looking strings
# (the ID, day, month, year, reminder) are not auto-converted to numbers.
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "3310658602400223725"
$ws.Cells.Item($row, 2).Value = "21"
$ws.Cells.Item($row, 3).Value = "10"
$ws.Cells.Item($row, 4).Value = "2022"
$ws.Cells.Item($row, 5).Value = "Nach Dänemark fahren"
$ws.Cells.Item($row, 6).Value = "1"

# The other data rows in this sheet carry no explicit cell style, so drop
# the temporary text format again now that the values are safely stored as
# text -- this keeps row 11 formatted the same as rows 3-10.
$rng.ClearFormats()
